# aggiornamento fino a 20/09/2021
# Append rows 375:385 with new daily data to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(375, 44449, 1, 11, 72.57851675903932),
    @(376, 44450, 0, 8, 52.78437582475588),
    @(377, 44451, 2, 7, 46.18632884666139),
    @(378, 44452, 0, 7, 46.18632884666139),
    @(379, 44453, 1, 5, 32.99023489047242),
    @(380, 44454, 0, 5, 32.99023489047242),
    @(381, 44455, 0, 4, 26.39218791237794),
    @(382, 44456, 5, 8, 52.78437582475588),
    @(383, 44457, 4, 12, 79.1765637371338),
    @(384, 44458, 2, 12, 79.1765637371338),
    @(385, 44459, 1, 13, 85.77461071522829)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}

# Copy the formatting (style) of the last existing row's date cell (A374)
# onto the newly added date cells (A375:A385), matching the s="2" style.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
$excel.CutCopyMode = $false
